$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.75
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 1.85
$ws.Range("J2").Value = 4.33
$ws.Range("K2").Value = 2.25
$ws.Range("L2").Value = 2.5
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 2.05
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 41
$ws.Range("AA2").Value = 29
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 15
$ws.Range("AI2").Value = 9.5
$ws.Range("AJ2").Value = 8.5
$ws.Range("AK2").Value = 15
$ws.Range("AL2").Value = 15
$ws.Range("AM2").Value = 23
$ws.Range("AN2").Value = 6
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 26
$ws.Range("AQ2").Value = 67
$ws.Range("AR2").Value = 81
$ws.Range("AW2").Value = 4
$ws.Range("AX2").Value = 10
$ws.Range("AY2").Value = 19
$ws.Range("AZ2").Value = 34

# Row 3
$ws.Range("G3").Value = 1.3
$ws.Range("H3").Value = 5.25
$ws.Range("I3").Value = 9.5
$ws.Range("J3").Value = 1.8
$ws.Range("K3").Value = 2.63
$ws.Range("L3").Value = 7.5
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 17
$ws.Range("O3").Value = 1.17
$ws.Range("P3").Value = 5
$ws.Range("Q3").Value = 1.57
$ws.Range("R3").Value = 2.35
$ws.Range("S3").Value = 1.29
$ws.Range("T3").Value = 3.5
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 8
$ws.Range("X3").Value = 7
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 8.5
$ws.Range("AA3").Value = 11
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 15
$ws.Range("AD3").Value = 10
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 301
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 41
$ws.Range("AJ3").Value = 23
$ws.Range("AK3").Value = 101
$ws.Range("AL3").Value = 51
$ws.Range("AM3").Value = 51
$ws.Range("AN3").Value = 3.4
$ws.Range("AO3").Value = 6
$ws.Range("AP3").Value = 17
$ws.Range("AQ3").Value = 15
$ws.Range("AR3").Value = 41
$ws.Range("AS3").Value = 101
$ws.Range("AT3").Value = 3.5
$ws.Range("AU3").Value = 9
$ws.Range("AV3").Value = 51
$ws.Range("AW3").Value = 9.5
$ws.Range("AX3").Value = 41
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 151
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 301

# Row 4
$ws.Range("G4").Value = 2.2
$ws.Range("H4").Value = 3.35
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 2.77
$ws.Range("K4").Value = 2.12
$ws.Range("L4").Value = 3.45
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 7.9
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.25
$ws.Range("Q4").Value = 1.75
$ws.Range("R4").Value = 1.87
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 2.99
$ws.Range("U4").Value = 1.6
$ws.Range("V4").Value = 2.05
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 11.25
$ws.Range("AA4").Value = 17
$ws.Range("AB4").Value = 25
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 12.5
$ws.Range("AF4").Value = 50
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 10.5
$ws.Range("AI4").Value = 16.5
$ws.Range("AL4").Value = 24
$ws.Range("AM4").Value = 29
$ws.Range("AN4").Value = 4.2
$ws.Range("AO4").Value = 11.25
$ws.Range("AP4").Value = 18.5
$ws.Range("AQ4").Value = 45
$ws.Range("AR4").Value = 70
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 2.8
$ws.Range("AU4").Value = 6.7
$ws.Range("AV4").Value = 55
$ws.Range("AW4").Value = 5
$ws.Range("AX4").Value = 15.5
$ws.Range("AY4").Value = 21
$ws.Range("AZ4").Value = 70
$ws.Range("BA4").Value = 100
$ws.Range("BB4").Value = 250

# Row 6
$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 2.1
$ws.Range("X6").Value = 12
$ws.Range("Y6").Value = 9.5
$ws.Range("AA6").Value = 19
$ws.Range("AY6").Value = 23

# Row 7
$ws.Range("G7").Value = 2.88
$ws.Range("I7").Value = 2.5
$ws.Range("J7").Value = 3.75
$ws.Range("L7").Value = 3.4
$ws.Range("W7").Value = 7
$ws.Range("X7").Value = 13
$ws.Range("AA7").Value = 29
$ws.Range("AD7").Value = 6
$ws.Range("AJ7").Value = 10
$ws.Range("AL7").Value = 23
$ws.Range("AO7").Value = 19
$ws.Range("AQ7").Value = 67
$ws.Range("AW7").Value = 4.33
$ws.Range("AX7").Value = 15

# Row 8
$ws.Range("G8").Value = 1.95
$ws.Range("I8").Value = 3.6
$ws.Range("J8").Value = 2.63
$ws.Range("L8").Value = 4.5
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 9.5
$ws.Range("X8").Value = 8.5
$ws.Range("AC8").Value = 9.5
$ws.Range("AG8").Value = 401
$ws.Range("AH8").Value = 9
$ws.Range("AJ8").Value = 12
$ws.Range("AR8").Value = 51
$ws.Range("AX8").Value = 23

# Row 10
$ws.Range("O10").Value = 1.22
$ws.Range("P10").Value = 4
$ws.Range("Q10").Value = 1.73
$ws.Range("R10").Value = 2.08

# Row 12
$ws.Range("G12").Value = 3.3
$ws.Range("H12").Value = 3.05
$ws.Range("I12").Value = 2.18
$ws.Range("J12").Value = 3.8
$ws.Range("K12").Value = 2.05
$ws.Range("L12").Value = 2.7
$ws.Range("M12").Value = 1.01
$ws.Range("N12").Value = 8.1
$ws.Range("U12").Value = 1.75
$ws.Range("V12").Value = 1.87
$ws.Range("W12").Value = 9
$ws.Range("X12").Value = 17
$ws.Range("Y12").Value = 11.5
$ws.Range("Z12").Value = 45
$ws.Range("AA12").Value = 32
$ws.Range("AD12").Value = 6
$ws.Range("AE12").Value = 14
$ws.Range("AH12").Value = 7.2
$ws.Range("AI12").Value = 10.5
$ws.Range("AJ12").Value = 8.75
$ws.Range("AK12").Value = 21
$ws.Range("AL12").Value = 18
$ws.Range("AN12").Value = 5.2
$ws.Range("AO12").Value = 18
$ws.Range("AP12").Value = 24
$ws.Range("AQ12").Value = 90
$ws.Range("AT12").Value = 2.55
$ws.Range("AU12").Value = 6.7
$ws.Range("AW12").Value = 4.05
$ws.Range("AX12").Value = 11
$ws.Range("AY12").Value = 18
$ws.Range("AZ12").Value = 40
$ws.Range("BA12").Value = 70
